$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns: J1 = "Kurum", K1 = "Üst Birim"
$ws.Cells.Item(1, 10).Value = "Kurum"
$ws.Cells.Item(1, 11).Value = "Üst Birim"

# Copy header style (s=2) from I1 to J1:K1
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats

# Column widths per diff
$ws.Columns.Item(10).ColumnWidth = 13.85546875
$ws.Columns.Item(11).ColumnWidth = 18

# Select K1 like the sheetView selection in the diff
$ws.Range("K1").Select()
